# Add "Insights" notes block below the Conforto/Carisma tables on Planilha1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("D18").Value = "Insights"
$ws.Range("E19").Value = "ta tudo igual (????????)"
$ws.Range("E20").Value = "Os participantes acharam a personagem confortável/normal"
$ws.Range("E21").Value = "Os participantes não acharam a personagem desconfortável"

# Match the saved selection/active cell state from the edited workbook.
$ws.Range("J16").Select()
